$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6: "Update CYRS Document" - Expected Delivery Date shifts from 2/4/2020 to 2/5/2020
$ws.Range("E6").Value = 43866

# Row 7: "Update HSI Document" - Expected Delivery Date shifts from 2/4/2020 to 2/5/2020
$ws.Range("E7").Value = 43866

# Row 8: "Update SRS Document" - Start Date shifts from 2/4/2020 to 2/5/2020
$ws.Range("D8").Value = 43866

# Row 8: Estimated Duration changes from "3 days" to "2 days"
$ws.Range("F8").Value = "2 days"

# Row 9: "Update RTM" - Review Status set to "Pending"
$ws.Range("J9").Value = "Pending"
